$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the window tab-ratio split slightly (299 -> 300 in the saved ratio).
$win = $wb.Windows.Item(1)
$win.TabRatio = 0.3

# Normalize the tab color (00FFFFFF -> FFFFFFFF, i.e. fully-opaque white).
$ws.Tab.Color = 16777215

# Add a new row with the label for the "7*row+col" table below the existing grid.
$ws.Range("B10").Value = "7*row+col"

# Move the active selection down to B11 (one row below the new content).
$ws.Range("B11").Select()
